$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "FAPs"
$ws.Range("B2").Value = "Lgi3"
$ws.Range("C2").Value = "Stx1a"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 3.0
$ws.Range("F2").Value = 1.0
$ws.Range("G2").Value = 0.6167776666666667
$ws.Range("H2").Value = 1.850333
$ws.Range("I2").Value = 0.7836323164322263
$ws.Range("J2").Value = 0.7836323164322262
$ws.Range("K2").Value = 3.0
$ws.Range("L2").Value = 1.0
$ws.Range("M2").Value = 2.618118666666667
$ws.Range("N2").Value = 7.854356
$ws.Range("O2").Value = 0.3278945548258764
$ws.Range("P2").Value = 0.3278945548258765
$ws.Range("Q2").Value = 1.614797122283111
$ws.Range("R2").Value = 14.533174100548
$ws.Range("S2").Value = 0.2569487695437151
$ws.Range("T2").Value = 0.2569487695437152

# Row 3
$ws.Range("A3").Value = "FAPs"
$ws.Range("B3").Value = "Lgi3"
$ws.Range("C3").Value = "Stx1a"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 3.0
$ws.Range("F3").Value = 1.0
$ws.Range("G3").Value = 0.6167776666666667
$ws.Range("H3").Value = 1.850333
$ws.Range("I3").Value = 0.7836323164322263
$ws.Range("J3").Value = 0.7836323164322262
$ws.Range("K3").Value = 3.0
$ws.Range("L3").Value = 1.0
$ws.Range("M3").Value = 4.293075666666667
$ws.Range("N3").Value = 12.879227
$ws.Range("O3").Value = 0.537667047898823
$ws.Range("P3").Value = 0.537667047898823
$ws.Range("Q3").Value = 2.647873192510111
$ws.Range("R3").Value = 23.830858732591
$ws.Range("S3").Value = 0.4213332742142314
$ws.Range("T3").Value = 0.4213332742142314

# Row 4
$ws.Range("A4").Value = "FAPs"
$ws.Range("B4").Value = "Lgi3"
$ws.Range("C4").Value = "Stx1a"
$ws.Range("D4").Value = "M2"
$ws.Range("E4").Value = 3.0
$ws.Range("F4").Value = 1.0
$ws.Range("G4").Value = 0.6167776666666667
$ws.Range("H4").Value = 1.850333
$ws.Range("I4").Value = 0.7836323164322263
$ws.Range("J4").Value = 0.7836323164322262
$ws.Range("K4").Value = 1.0
$ws.Range("L4").Value = 0.3333333333333333
$ws.Range("M4").Value = 0.01260166666666667
$ws.Range("N4").Value = 0.037805
$ws.Range("O4").Value = 0.001578239341989624
$ws.Range("P4").Value = 0.001578239341989624
$ws.Range("Q4").Value = 0.007772426562777777
$ws.Range("R4").Value = 0.069951839065
$ws.Range("S4").Value = 0.001236759351447802
$ws.Range("T4").Value = 0.001236759351447802

# Row 5
$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Lgi3"
$ws.Range("C5").Value = "Stx1a"
$ws.Range("D5").Value = "sCs"
$ws.Range("E5").Value = 3.0
$ws.Range("F5").Value = 1.0
$ws.Range("G5").Value = 0.6167776666666667
$ws.Range("H5").Value = 1.850333
$ws.Range("I5").Value = 0.7836323164322263
$ws.Range("J5").Value = 0.7836323164322262
$ws.Range("K5").Value = 3.0
$ws.Range("L5").Value = 1.0
$ws.Range("M5").Value = 1.06084
$ws.Range("N5").Value = 3.18252
$ws.Range("O5").Value = 0.1328601579333109
$ws.Range("P5").Value = 0.1328601579333109
$ws.Range("Q5").Value = 0.6543024199066667
$ws.Range("R5").Value = 5.88872177916
$ws.Range("S5").Value = 0.1041135133228319
$ws.Range("T5").Value = 0.1041135133228319

# Row 6
$ws.Range("A6").Value = "sCs"
$ws.Range("B6").Value = "Lgi3"
$ws.Range("C6").Value = "Stx1a"
$ws.Range("D6").Value = "ECs"
$ws.Range("E6").Value = 2.0
$ws.Range("F6").Value = 0.6666666666666666
$ws.Range("G6").Value = 0.1702976666666667
$ws.Range("H6").Value = 0.510893
$ws.Range("I6").Value = 0.2163676835677737
$ws.Range("J6").Value = 0.2163676835677737
$ws.Range("K6").Value = 3.0
$ws.Range("L6").Value = 1.0
$ws.Range("M6").Value = 2.618118666666667
$ws.Range("N6").Value = 7.854356
$ws.Range("O6").Value = 0.3278945548258764
$ws.Range("P6").Value = 0.3278945548258765
$ws.Range("Q6").Value = 0.4458594999897778
$ws.Range("R6").Value = 4.012735499908
$ws.Range("S6").Value = 0.07094578528216125
$ws.Range("T6").Value = 0.07094578528216126

# Row 7
$ws.Range("A7").Value = "sCs"
$ws.Range("B7").Value = "Lgi3"
$ws.Range("C7").Value = "Stx1a"
$ws.Range("D7").Value = "FAPs"
$ws.Range("E7").Value = 2.0
$ws.Range("F7").Value = 0.6666666666666666
$ws.Range("G7").Value = 0.1702976666666667
$ws.Range("H7").Value = 0.510893
$ws.Range("I7").Value = 0.2163676835677737
$ws.Range("J7").Value = 0.2163676835677737
$ws.Range("K7").Value = 3.0
$ws.Range("L7").Value = 1.0
$ws.Range("M7").Value = 4.293075666666667
$ws.Range("N7").Value = 12.879227
$ws.Range("O7").Value = 0.537667047898823
$ws.Range("P7").Value = 0.537667047898823
$ws.Range("Q7").Value = 0.7311007688567779
$ws.Range("R7").Value = 6.579906919711001
$ws.Range("S7").Value = 0.1163337736845916
$ws.Range("T7").Value = 0.1163337736845916

# Row 8
$ws.Range("A8").Value = "sCs"
$ws.Range("B8").Value = "Lgi3"
$ws.Range("C8").Value = "Stx1a"
$ws.Range("D8").Value = "M2"
$ws.Range("E8").Value = 2.0
$ws.Range("F8").Value = 0.6666666666666666
$ws.Range("G8").Value = 0.1702976666666667
$ws.Range("H8").Value = 0.510893
$ws.Range("I8").Value = 0.2163676835677737
$ws.Range("J8").Value = 0.2163676835677737
$ws.Range("K8").Value = 1.0
$ws.Range("L8").Value = 0.3333333333333333
$ws.Range("M8").Value = 0.01260166666666667
$ws.Range("N8").Value = 0.037805
$ws.Range("O8").Value = 0.001578239341989624
$ws.Range("P8").Value = 0.001578239341989624
$ws.Range("Q8").Value = 0.002146034429444445
$ws.Range("R8").Value = 0.019314309865
$ws.Range("S8").Value = 0.0003414799905418224
$ws.Range("T8").Value = 0.0003414799905418224

# Row 9
$ws.Range("A9").Value = "sCs"
$ws.Range("B9").Value = "Lgi3"
$ws.Range("C9").Value = "Stx1a"
$ws.Range("D9").Value = "sCs"
$ws.Range("E9").Value = 2.0
$ws.Range("F9").Value = 0.6666666666666666
$ws.Range("G9").Value = 0.1702976666666667
$ws.Range("H9").Value = 0.510893
$ws.Range("I9").Value = 0.2163676835677737
$ws.Range("J9").Value = 0.2163676835677737
$ws.Range("K9").Value = 3.0
$ws.Range("L9").Value = 1.0
$ws.Range("M9").Value = 1.06084
$ws.Range("N9").Value = 3.18252
$ws.Range("O9").Value = 0.1328601579333109
$ws.Range("P9").Value = 0.1328601579333109
$ws.Range("Q9").Value = 0.1806585767066667
$ws.Range("R9").Value = 1.62592719036
$ws.Range("S9").Value = 0.02874664461047906
$ws.Range("T9").Value = 0.02874664461047906
